{"js": "// Office.js (Word JavaScript API) edit matching the target diff:\n//  - \"Trocar icones.\" becomes a blank, unstyled paragraph (list formatting removed).\n//  - \"Ajustar responsividade de imagens \"check\" no site.\" is removed entirely.\n//  - \"Telas de login e cadastro para site e mobile\" is removed entirely.\n//  - \"Area do cliente no site e no mobile.\" is removed entirely.\n//  - \"Area de servi\u00e7os website\" is removed entirely.\n//  - \"Alinhar texto \"quem somos\" no website\" is removed entirely.\n//  - \"Tela de Pedidos\" is removed entirely.\n//  - \"Tela de Carrinho\" is removed entirely.\n// All other paragraphs (e.g. \"Adicionar o whatsapp...\", \"Trocar t\u00edtulo...\",\n// \"Shared preferences e session\", \"Emails de promo\u00e7\u00f5es...\", \"Async e\n// seguran\u00e7a\", \"Tela de Ajuda(...)\") keep the same visible text \u2014 Word's\n// proofing engine is the only thing that later sprinkles <w:proofErr/> spell\n// check markers and extra run splits around them, which carries no visible\n// or textual change, so we leave that text untouched here.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Paragraphs whose text must disappear completely (paragraph + its text).\nconst textsToRemove = new Set([\n  \"Ajustar responsividade de imagens \\u201ccheck\\u201d no site.\",\n  \"Telas de login e cadastro para site e mobile\",\n  \"Area do cliente no site e no mobile.\",\n  \"Area de servi\\u00e7os website\",\n  \"Alinhar texto \\u201cquem somos\\u201d no website\",\n  \"Tela de Pedidos\",\n  \"Tela de Carrinho\",\n]);\n\n// The very first paragraph (\"Trocar icones.\") turns into an empty paragraph\n// with no list numbering / paragraph style (equivalent to a bare <w:p/>).\nconst firstParaText = \"Trocar icones.\";\n\nconst toDelete = [];\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === firstParaText) {\n    paragraph.clear();\n    paragraph.style = \"Normal\";\n  } else if (textsToRemove.has(paragraph.text)) {\n    toDelete.push(paragraph);\n  }\n}\n\nfor (const paragraph of toDelete) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit matching the target diff:\n#  - \"Trocar icones.\" becomes a blank, unstyled paragraph (list formatting\n#    removed) - equivalent to an empty <w:p/>.\n#  - \"Ajustar responsividade de imagens \"check\" no site.\" is removed entirely.\n#  - \"Telas de login e cadastro para site e mobile\" is removed entirely.\n#  - \"Area do cliente no site e no mobile.\" is removed entirely.\n#  - \"Area de servi\u00e7os website\" is removed entirely.\n#  - \"Alinhar texto \"quem somos\" no website\" is removed entirely.\n#  - \"Tela de Pedidos\" is removed entirely.\n#  - \"Tela de Carrinho\" is removed entirely.\n# All other paragraphs (e.g. \"Adicionar o whatsapp...\", \"Trocar t\u00edtulo...\",\n# \"Shared preferences e session\", \"Emails de promo\u00e7\u00f5es...\", \"Async e\n# seguran\u00e7a\", \"Tela de Ajuda(...)\") keep the same visible text - Word's\n# proofing engine is the only thing that later sprinkles <w:proofErr/> spell\n# check markers and extra run splits around them, which carries no visible\n# or textual change, so that text is left untouched here.\n\n$d = $word.ActiveDocument\n\n$textsToRemove = @(\n    \"Ajustar responsividade de imagens \u201ccheck\u201d no site.\",\n    \"Telas de login e cadastro para site e mobile\",\n    \"Area do cliente no site e no mobile.\",\n    \"Area de servi\u00e7os website\",\n    \"Alinhar texto \u201cquem somos\u201d no website\",\n    \"Tela de Pedidos\",\n    \"Tela de Carrinho\"\n)\n\n$firstParaText = \"Trocar icones.\"\n\n# Walk paragraphs back-to-front so deleting one never shifts the index of a\n# paragraph we still need to visit.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $firstParaText) {\n        $p.Range.Text = \"\"\n        $p.Range.Style = \"Normal\"\n    } elseif ($textsToRemove -contains $t) {\n        $p.Range.Delete()\n    }\n}\n"}
